$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for "Empresa"
$ws.Range("A1").EntireColumn.Insert()

# Apply the existing header style (bold/border/centered) to the new A1 and F1
# header cells by copying the format from the already-styled B1 header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row text
$ws.Range("A1").Value = "Empresa"
$ws.Range("F1").Value = "Classificação"

# Data
$data = @(
    @("Shein", "Respondeu 99.9% das reclamações recebidas.", "Dos que avaliaram, 62.8% voltariam a fazer negócio.", "A empresa resolveu 67.7% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como BOM. A nota média nos últimos 6 meses é 7.0/10.", "1º Melhor"),
    @("Centauro", "Respondeu 100% das reclamações recebidas.", "Dos que avaliaram, 61% voltariam a fazer negócio.", "A empresa resolveu 82% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como BOM. A nota média nos últimos 6 meses é 7.4/10.", "2º Melhor"),
    @("Lojas Renner", "Respondeu 99.9% das reclamações recebidas.", "Dos que avaliaram, 80% voltariam a fazer negócio.", "A empresa resolveu 89.6% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como ÓTIMO. A nota média nos últimos 6 meses é 8.5/10.", "3º Melhor"),
    @("Yeesco", "Respondeu 59.9% das reclamações recebidas.", "Dos que avaliaram, 9.1% voltariam a fazer negócio.", "A empresa resolveu 20.5% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como Não Recomendada. A nota média nos últimos 6 meses é 2.3/10.", "1º Pior"),
    @("bycih store", "Respondeu 99.8% das reclamações recebidas.", "Dos que avaliaram, 13.1% voltariam a fazer negócio.", "A empresa resolveu 37.3% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como Não Recomendada. A nota média nos últimos 6 meses é 3.9/10.", "2º Pior"),
    @("Miazzi", "Respondeu 99.7% das reclamações recebidas.", "Dos que avaliaram, 5.4% voltariam a fazer negócio.", "A empresa resolveu 12% das reclamações recebidas.", "O consumidor avaliou o atendimento dessa empresa como Não Recomendada. A nota média nos últimos 6 meses é 2.7/10.", "3º Pior")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}
